# Updates crypto Price (D) and Volume(1h) (E) columns with latest scraped values.
# Cells store text (not numeric) values, e.g. "323.08" and "8.66%", so we force
# a Text number format before assigning to avoid Excel auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '323.08'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '8.66%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '49.59'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '18.20%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.350'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '6.76%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08160'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '8.50%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.606'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '5.31%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.666'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '5.20%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.176'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '27.02%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1349'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '13.00%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1959'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '7.20%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09608'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '8.00%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04536'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '11.21%'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.04%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001325'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '3.49%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005909'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.22%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.398'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '1.26%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.435'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.40%'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '2.43%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.184'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1.10%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1419'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '4.32%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.3051'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-1.64%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04303'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '5.34%'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '3.13%'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '9.66%'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '9.70%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0003720'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.11%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02780'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '15.03%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05571'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '7.05%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.006296'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.12%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007686'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.20%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1449'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '9.11%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.007690'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '4.04%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008061'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '11.25%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3515'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '18.41%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006776'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '2.90%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.12%'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '93.61%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.003998'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-4.83%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002099'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.12%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0001999'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.12%'

$wb.Save()
